# Plague.xlsx edit: add additional advisors / event summary rows (15-19)
# Mirrors the OUTCOME row (row 6) styling for the new "header" cells in column A
# and the wrapped merged-text cells in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the formatting of the existing OUTCOME row (A6:M6) down onto the five
#    new rows so fonts / fills / borders / alignment match the rest of the sheet.
$ws.Range("A6:M6").Copy()
$ws.Range("A15:M19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Column A headers for the new block.
$ws.Range("A15").Value = "EVENT SUMMARY HEADER"
$ws.Range("A16").Value = "ACTION 1 EVENT SUMMARY"
$ws.Range("A17").Value = "ACTION 2 EVENT SUMMARY"
$ws.Range("A18").Value = "ACTION 3 EVENT SUMMARY"
$ws.Range("A19").Value = "ACTION 4 EVENT SUMMARY"

# 3) Column B body text for the new block.
$ws.Range("B15").Value = "You have received news that the disease has sweeped through the neighbouring kingdoms and situation has gotten severe."
$ws.Range("B16").Value = "Thankfully, the strict border restrictions meant that you were able to control the flow of people entering your Kingdom, minimising any risk of any potential foreign infection. The stock on the medicine also meant that your Kingdom is able to cope should an outbreak occur."
$ws.Range("B17").Value = "Although this had helped to reduce the number of infected people, such actions were viewed to be treacherous by the neighbouring Kingdoms, ultimately angering them."
$ws.Range("B18").Value = "The decision to cut off interactions with neighbouring Kingdoms meant that you removed the risk of any infected people from the neighbouring Kingdoms entering your own. However, this came a cost of your Kingdom's trade and relationship ties with the neighbouring Kingdoms."
$ws.Range("B19").Value = "Letting your guard down was not the best as you unknowningly let infected people from the neighbouring Kingdoms enter your own and allowed for spreading of the disease."

# 4) Merge the B:M span on each new row, same pattern as the other text rows.
$ws.Range("B15:M15").Merge()
$ws.Range("B16:M16").Merge()
$ws.Range("B17:M17").Merge()
$ws.Range("B18:M18").Merge()
$ws.Range("B19:M19").Merge()

# 5) Give the text rows a bit more room since the copy is longer than the
#    single-line OUTCOME banner it was copied from.
$ws.Range("15:15").RowHeight = 55.2
$ws.Range("16:16").RowHeight = 72.6
$ws.Range("17:17").RowHeight = 72.6
$ws.Range("18:18").RowHeight = 72.6
$ws.Range("19:19").RowHeight = 72.6

# 6) Scroll the view down to the newly added rows, matching the author's
#    final selection/viewport when they saved.
$ws.Range("B20").Select()
$excel.ActiveWindow.ScrollRow = 15
